# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# This script updates the "Estado de Cuenta" sheet for NIT 9005464341:
#  - Removes the BERENICE BRAVO BLANCO entries (periods 1809-1812) entirely.
#  - Moves ALFONSO ENRIQUE RAMOS DIAZ's single entry (period 1709) to the top
#    of the worker table.
#  - Keeps MARLYNG VELEZ BLANQUICETT's 13 periods, but re-orders them from
#    ascending 1901 -> 1912 -> 2001 and refreshes the "Valor Mora" (F) and
#    "Salario Basico" (G) amounts.
#  - Refreshes the summary totals (Valor Mora total, worker count, period
#    count) to reflect the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the 5 rows belonging to BERENICE BRAVO BLANCO (old rows 29-33).
#    This automatically shifts the ALFONSO row (old row 34) up to row 29,
#    and shifts the footer rows (old 39/40) up to 34/35, along with the
#    dimension / merged-cell ranges.
# ---------------------------------------------------------------------
$ws.Rows("29:33").Delete()

# ---------------------------------------------------------------------
# 2) Row 16 now becomes ALFONSO ENRIQUE RAMOS DIAZ's single entry
#    (moved from the bottom of the table to the top).
# ---------------------------------------------------------------------
$ws.Range("B16").Value() = "CC"
$ws.Range("C16").Value() = "7921211"
$ws.Range("D16").Value() = "ALFONSO ENRIQUE RAMOS DIAZ"
$ws.Range("E16").Value() = "1709"
$ws.Range("F16").Value() = 29509
$ws.Range("G16").Value() = 737717

# ---------------------------------------------------------------------
# 3) Rows 17-29 hold MARLYNG VELEZ BLANQUICETT's 13 periods, now sorted
#    ascending (1901 .. 1912, 2001) with refreshed Valor Mora / Salario
#    Basico figures.
# ---------------------------------------------------------------------
$marlyngPeriods = @("1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001")
$marlyngValorMora = @(11042,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125,33125)
$marlyngSalario = 877803

$row = 17
for ($i = 0; $i -lt $marlyngPeriods.Length; $i++) {
    $ws.Range("B$row").Value() = "CC"
    $ws.Range("C$row").Value() = "45515194"
    $ws.Range("D$row").Value() = "MARLYNG VELEZ BLANQUICETT"
    $ws.Range("E$row").Value() = $marlyngPeriods[$i]
    $ws.Range("F$row").Value() = $marlyngValorMora[$i]
    $ws.Range("G$row").Value() = $marlyngSalario
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 4) Refresh the summary figures at the top of the sheet:
#    - E11: total Valor Mora
#    - C13: Cant. Trabajadores (worker count)
#    - F13: Cant. Periodos (period count)
# ---------------------------------------------------------------------
$ws.Range("E11").Value() = 438051
$ws.Range("C13").Value() = 2
$ws.Range("F13").Value() = 14
